# The "建物" (building) sheet's property_category column (I) was mistakenly
# filled with "land" (copied from the 土地 sheet). Correct every data row so
# it reads "building" instead, matching the sheet it actually belongs to.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 9).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 17 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($cell.Value2 -eq "land") {
        $cell.Value2 = "building"
    }
}
